$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "lexion 670" row (row 5, the "C" clutch entry): narrow the clutch band
# from [0.1, 0.7] to [0, 0.8]. Recalculated (dependent) cells G5:P5 and the
# matching chart series update automatically.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0.8

# Restore the last recorded selection on the sheet.
$ws.Range("D6").Select()
